# Applies the cryptos.xlsx price/volume/rank update described by the commit:
# "Updated cryptos list on Mon Jul 17 05:27:21 UTC 2023 with GitHub Actions"
#
# Column D holds price strings and column E holds percent-change strings that
# are stored as literal text in the workbook (not real numbers/percentages).
# Several D values look like plain numbers (e.g. "1.000", "0.7494") so the
# COM layer would otherwise coerce them to numeric doubles and silently drop
# the trailing zeros / exact formatting. Pre-setting NumberFormat to "@" (Text)
# on those specific cells keeps the assignment a literal string, matching the
# source data exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.308.39'
$ws.Range('E2').Value = '  +0.26%  '

$ws.Range('D3').Value = '1.933.81'
$ws.Range('E3').Value = '  +0.33%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7494'
$ws.Range('E5').Value = '  +5.13%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.20'
$ws.Range('E6').Value = '  -2.21%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.19'
$ws.Range('E8').Value = '  +2.99%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3184'
$ws.Range('E9').Value = '  -0.66%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07037'
$ws.Range('E10').Value = '  -0.54%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7837'
$ws.Range('E11').Value = '  -1.04%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08038'

$ws.Range('D13').Value = '1.937.00'
$ws.Range('E13').Value = '  +0.40%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.403'
$ws.Range('E14').Value = '  +0.46%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '93.04'
$ws.Range('E15').Value = '  -1.92%  '

$ws.Range('E16').Value = '  -0.96%  '

$ws.Range('D17').Value = '30.317.57'
$ws.Range('E17').Value = '  +0.22%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.083'
$ws.Range('E18').Value = '  +5.61%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '252.50'
$ws.Range('E19').Value = '  -1.61%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007988'
$ws.Range('E20').Value = '  -0.47%  '

$ws.Range('D21').Value = '2.189.80'
$ws.Range('E21').Value = '  +0.40%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9989'
$ws.Range('E22').Value = '  -0.13%  '

$ws.Range('E23').Value = '  -0.06%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.711'
$ws.Range('E24').Value = '  -1.91%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.587'
$ws.Range('E25').Value = '  +0.64%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.91'
$ws.Range('E26').Value = '  -0.75%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.08'
$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1308'
$ws.Range('E28').Value = '  +3.60%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.210'
$ws.Range('E29').Value = '  -2.32%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.367'
$ws.Range('E30').Value = '  +0.78%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.546'
$ws.Range('E31').Value = '  +1.08%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.442'
$ws.Range('E32').Value = '  +1.12%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.146'
$ws.Range('E33').Value = '  +0.56%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05279'
$ws.Range('E34').Value = '  +2.64%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.334'
$ws.Range('E35').Value = '  +5.08%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7572'
$ws.Range('E36').Value = '  +1.63%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.782'
$ws.Range('E37').Value = '  +0.77%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01962'
$ws.Range('E38').Value = '  +0.22%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.805'
$ws.Range('E39').Value = '  +0.26%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '79.04'
$ws.Range('E40').Value = '  +1.99%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.534'
$ws.Range('E41').Value = '  +2.81%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4516'
$ws.Range('E42').Value = '  +0.53%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.985'
$ws.Range('E43').Value = '  -0.15%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8390'

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.10'
$ws.Range('E46').Value = '  +3.36%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.711'
$ws.Range('E47').Value = '  +3.72%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.69'
$ws.Range('E48').Value = '  +1.21%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.64'
$ws.Range('E49').Value = '  +3.15%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1252'
$ws.Range('E50').Value = '  +10.38%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '967.20'
$ws.Range('E51').Value = '  +5.89%  '
